# Apply the recorded edit: rows 9-16,18 of the "Artfynd" sheet have their
# species-observation data (columns A,B,D,E,F,G,H,Q,R) permuted between rows.
# Row 17 is untouched.
#
# Permutation discovered from the diff (before-row -> after-row, i.e. the
# values that used to live in row X now live in row Y):
#   9  -> 13
#   13 -> 9
#   10 -> 12
#   12 -> 14
#   14 -> 15
#   15 -> 10
#   11 -> 18
#   18 -> 16
#   16 -> 11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "before" values for the columns that change, for every
# affected row, so we can write them into their destination row without
# clobbering data we still need to read.
$rows = @(9, 10, 11, 12, 13, 14, 15, 16, 18)

$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = @{
        A = $ws.Cells.Item($r, 1).Value2
        B = $ws.Cells.Item($r, 2).Value2
        D = $ws.Cells.Item($r, 4).Value2
        E = $ws.Cells.Item($r, 5).Value2
        F = $ws.Cells.Item($r, 6).Value2
        G = $ws.Cells.Item($r, 7).Value2
        H = $ws.Cells.Item($r, 8).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        R = $ws.Cells.Item($r, 18).Value2
    }
}

# before-row -> after-row mapping
$mapping = @{
    9  = 13
    10 = 12
    11 = 18
    12 = 14
    13 = 9
    14 = 15
    15 = 10
    16 = 11
    18 = 16
}

foreach ($src in $mapping.Keys) {
    $dst = $mapping[$src]
    $data = $snapshot[$src]

    $ws.Cells.Item($dst, 1).Value = $data.A
    $ws.Cells.Item($dst, 2).Value = $data.B
    $ws.Cells.Item($dst, 4).Value = $data.D
    $ws.Cells.Item($dst, 5).Value = $data.E
    $ws.Cells.Item($dst, 6).Value = $data.F
    $ws.Cells.Item($dst, 7).Value = $data.G
    $ws.Cells.Item($dst, 8).Value = $data.H
    $ws.Cells.Item($dst, 17).Value = $data.Q
    $ws.Cells.Item($dst, 18).Value = $data.R
}

# Column L (12) is an otherwise-unused helper column that is blank on every
# affected row, except row 14 originally had no L cell at all while row 12
# had a present-but-empty one. Since row 12's data moves into row 14, and
# row 14's (L-less) data moves into row 15, the net effect is: row 14 gains
# a blank L cell and row 15 ends up with no L cell at all.
$ws.Cells.Item(14, 12).Value = ""
$ws.Cells.Item(15, 12).ClearContents()
